# =====================================================================
# Add 2022-Q3 data
#  1. Insert a new worksheet named "2022-Q3" right before the existing
#     "2022-Q2" sheet, and fill it with the fund holdings table.
#  2. On the summary sheet ("总计") insert a new row for 2022-Q3 at the
#     top of the data (row 2), shifting the existing quarters down.
# =====================================================================

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)      # "总计"
$q2Sheet = $wb.Worksheets.Item(2)      # "2022-Q2" (currently 2nd tab)

# ---------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet before the "2022-Q2" sheet
# ---------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# Bring over the header style (bold font + border) used on every
# other quarter sheet's header row, and the index-column style used
# for column A, by copying formats from the summary sheet.
$summary.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

# ---- Header row (B1:H1), text with bold/border style copied from summary sheet ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- Data rows 2-16 ----
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'005583"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "'易方达港股通红利灵活配置混合"
$newSheet.Range("C2").Style = "Normal"
$newSheet.Range("D2").Value = "'6.98"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'90.31"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'4.18"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.2918"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'013991"
$newSheet.Range("B3").Style = "Normal"
$newSheet.Range("C3").Value = "'中欧港股通精选一年持有混合A"
$newSheet.Range("C3").Style = "Normal"
$newSheet.Range("D3").Value = "'6.69"
$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("E3").Value = "'93.38"
$newSheet.Range("E3").Style = "Normal"
$newSheet.Range("F3").Value = "'3.54"
$newSheet.Range("F3").Style = "Normal"
$newSheet.Range("G3").Value = "'0.2368"
$newSheet.Range("G3").Style = "Normal"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'007592"
$newSheet.Range("B4").Style = "Normal"
$newSheet.Range("C4").Value = "'华夏价值精选混合"
$newSheet.Range("C4").Style = "Normal"
$newSheet.Range("D4").Value = "'2.26"
$newSheet.Range("D4").Style = "Normal"
$newSheet.Range("E4").Value = "'93.77"
$newSheet.Range("E4").Style = "Normal"
$newSheet.Range("F4").Value = "'10.29"
$newSheet.Range("F4").Style = "Normal"
$newSheet.Range("G4").Value = "'0.2326"
$newSheet.Range("G4").Style = "Normal"
$newSheet.Range("H4").Value = 1

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'006049"
$newSheet.Range("B5").Style = "Normal"
$newSheet.Range("C5").Value = "'恒越研究精选混合A/B"
$newSheet.Range("C5").Style = "Normal"
$newSheet.Range("D5").Value = "'5.49"
$newSheet.Range("D5").Style = "Normal"
$newSheet.Range("E5").Value = "'89.87"
$newSheet.Range("E5").Style = "Normal"
$newSheet.Range("F5").Value = "'3.64"
$newSheet.Range("F5").Style = "Normal"
$newSheet.Range("G5").Value = "'0.1998"
$newSheet.Range("G5").Style = "Normal"
$newSheet.Range("H5").Value = 6

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'013992"
$newSheet.Range("B6").Style = "Normal"
$newSheet.Range("C6").Value = "'中欧港股通精选一年持有混合C"
$newSheet.Range("C6").Style = "Normal"
$newSheet.Range("D6").Value = "'4.68"
$newSheet.Range("D6").Style = "Normal"
$newSheet.Range("E6").Value = "'93.38"
$newSheet.Range("E6").Style = "Normal"
$newSheet.Range("F6").Value = "'3.54"
$newSheet.Range("F6").Style = "Normal"
$newSheet.Range("G6").Value = "'0.1657"
$newSheet.Range("G6").Style = "Normal"
$newSheet.Range("H6").Value = 9

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'001581"
$newSheet.Range("B7").Style = "Normal"
$newSheet.Range("C7").Value = "'华安沪港深通精选混合A"
$newSheet.Range("C7").Style = "Normal"
$newSheet.Range("D7").Value = "'4.87"
$newSheet.Range("D7").Style = "Normal"
$newSheet.Range("E7").Value = "'87.45"
$newSheet.Range("E7").Style = "Normal"
$newSheet.Range("F7").Value = "'2.92"
$newSheet.Range("F7").Style = "Normal"
$newSheet.Range("G7").Value = "'0.1422"
$newSheet.Range("G7").Style = "Normal"
$newSheet.Range("H7").Value = 10

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'007192"
$newSheet.Range("B8").Style = "Normal"
$newSheet.Range("C8").Value = "'恒越研究精选混合C"
$newSheet.Range("C8").Style = "Normal"
$newSheet.Range("D8").Value = "'3.73"
$newSheet.Range("D8").Style = "Normal"
$newSheet.Range("E8").Value = "'89.87"
$newSheet.Range("E8").Style = "Normal"
$newSheet.Range("F8").Value = "'3.64"
$newSheet.Range("F8").Style = "Normal"
$newSheet.Range("G8").Value = "'0.1358"
$newSheet.Range("G8").Style = "Normal"
$newSheet.Range("H8").Value = 6

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'160125"
$newSheet.Range("B9").Style = "Normal"
$newSheet.Range("C9").Value = "'南方香港优选股票（QDII-LOF）"
$newSheet.Range("C9").Style = "Normal"
$newSheet.Range("D9").Value = "'2.01"
$newSheet.Range("D9").Style = "Normal"
$newSheet.Range("E9").Value = "'81.74"
$newSheet.Range("E9").Style = "Normal"
$newSheet.Range("F9").Value = "'3.34"
$newSheet.Range("F9").Style = "Normal"
$newSheet.Range("G9").Value = "'0.0671"
$newSheet.Range("G9").Style = "Normal"
$newSheet.Range("H9").Value = 6

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'012993"
$newSheet.Range("B10").Style = "Normal"
$newSheet.Range("C10").Value = "'汇添富品牌力一年持有混合A"
$newSheet.Range("C10").Style = "Normal"
$newSheet.Range("D10").Value = "'1.78"
$newSheet.Range("D10").Style = "Normal"
$newSheet.Range("E10").Value = "'64.64"
$newSheet.Range("E10").Style = "Normal"
$newSheet.Range("F10").Value = "'3.64"
$newSheet.Range("F10").Style = "Normal"
$newSheet.Range("G10").Value = "'0.0648"
$newSheet.Range("G10").Style = "Normal"
$newSheet.Range("H10").Value = 3

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "'519601"
$newSheet.Range("B11").Style = "Normal"
$newSheet.Range("C11").Value = "'海富通中国海外精选混合（QDII）"
$newSheet.Range("C11").Style = "Normal"
$newSheet.Range("D11").Value = "'0.51"
$newSheet.Range("D11").Style = "Normal"
$newSheet.Range("E11").Value = "'73.52"
$newSheet.Range("E11").Style = "Normal"
$newSheet.Range("F11").Value = "'3.99"
$newSheet.Range("F11").Style = "Normal"
$newSheet.Range("G11").Value = "'0.0203"
$newSheet.Range("G11").Style = "Normal"
$newSheet.Range("H11").Value = 6

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "'007518"
$newSheet.Range("B12").Style = "Normal"
$newSheet.Range("C12").Value = "'东方阿尔法优选混合A"
$newSheet.Range("C12").Style = "Normal"
$newSheet.Range("D12").Value = "'1.04"
$newSheet.Range("D12").Style = "Normal"
$newSheet.Range("E12").Value = "'84.81"
$newSheet.Range("E12").Style = "Normal"
$newSheet.Range("F12").Value = "'1.64"
$newSheet.Range("F12").Style = "Normal"
$newSheet.Range("G12").Value = "'0.0171"
$newSheet.Range("G12").Style = "Normal"
$newSheet.Range("H12").Value = 5

$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "'012994"
$newSheet.Range("B13").Style = "Normal"
$newSheet.Range("C13").Value = "'汇添富品牌力一年持有混合C"
$newSheet.Range("C13").Style = "Normal"
$newSheet.Range("D13").Value = "'0.26"
$newSheet.Range("D13").Style = "Normal"
$newSheet.Range("E13").Value = "'64.64"
$newSheet.Range("E13").Style = "Normal"
$newSheet.Range("F13").Value = "'3.64"
$newSheet.Range("F13").Style = "Normal"
$newSheet.Range("G13").Value = "'0.0095"
$newSheet.Range("G13").Style = "Normal"
$newSheet.Range("H13").Value = 3

$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "'007519"
$newSheet.Range("B14").Style = "Normal"
$newSheet.Range("C14").Value = "'东方阿尔法优选混合C"
$newSheet.Range("C14").Style = "Normal"
$newSheet.Range("D14").Value = "'0.47"
$newSheet.Range("D14").Style = "Normal"
$newSheet.Range("E14").Value = "'84.81"
$newSheet.Range("E14").Style = "Normal"
$newSheet.Range("F14").Value = "'1.64"
$newSheet.Range("F14").Style = "Normal"
$newSheet.Range("G14").Value = "'0.0077"
$newSheet.Range("G14").Style = "Normal"
$newSheet.Range("H14").Value = 5

$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "'519602"
$newSheet.Range("B15").Style = "Normal"
$newSheet.Range("C15").Value = "'海富通大中华精选混合（QDII）"
$newSheet.Range("C15").Style = "Normal"
$newSheet.Range("D15").Value = "'0.10"
$newSheet.Range("D15").Style = "Normal"
$newSheet.Range("E15").Value = "'87.37"
$newSheet.Range("E15").Style = "Normal"
$newSheet.Range("F15").Value = "'5.19"
$newSheet.Range("F15").Style = "Normal"
$newSheet.Range("G15").Value = "'0.0052"
$newSheet.Range("G15").Style = "Normal"
$newSheet.Range("H15").Value = 4

$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "'016289"
$newSheet.Range("B16").Style = "Normal"
$newSheet.Range("C16").Value = "'华安沪港深通精选混合C"
$newSheet.Range("C16").Style = "Normal"
$newSheet.Range("D16").Value = "'0.00"
$newSheet.Range("D16").Style = "Normal"
$newSheet.Range("E16").Value = "'87.45"
$newSheet.Range("E16").Style = "Normal"
$newSheet.Range("F16").Value = "'2.92"
$newSheet.Range("F16").Style = "Normal"
$newSheet.Range("G16").Value = 0
$newSheet.Range("H16").Value = 10

# ---------------------------------------------------------------
# 2. Update the summary sheet ("总计"): add a new row for 2022-Q3 at
#    the top of the data (row 2), pushing the existing quarters
#    (2022-Q2, 2021-Q4, 2021-Q3, 2020-Q4) down by one row each.
#    Shifting is done with Copy/PasteSpecial (bottom row first) so
#    every cell keeps exactly the style it already had, instead of
#    using Rows.Insert() (which would blend row 1's/row 2's borders
#    into a brand-new, unused style).
# ---------------------------------------------------------------

# Row 6 is brand new: give its index cell (column A) the same style
# used by the rest of the index column before copying any values.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)

# Shift rows 2-5 down into rows 3-6, starting from the bottom so we
# never overwrite a row before it has been copied.
$summary.Range("A5:D5").Copy()
$summary.Range("A6:D6").PasteSpecial(-4104)
$summary.Range("A4:D4").Copy()
$summary.Range("A5:D5").PasteSpecial(-4104)
$summary.Range("A3:D3").Copy()
$summary.Range("A4:D4").PasteSpecial(-4104)
$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4104)

# Write the new 2022-Q3 summary row into row 2 (style is already
# correct, carried over from the former row 2).
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 1.6
